$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "pen"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 0.99
